$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update notes/date for Task 9 (row 10)
$ws.Range("C10").Value = "regsiter backend complete, needs error message to refleft when misinput"

# D10 must stay a plain text string "02/14/2024" (not get auto-converted to a
# date serial). Build it as a text formula result, then paste-special as
# values so the stored cell is a literal shared string with no special
# number format applied (matches original unstyled text cell).
$ws.Range("F1").Formula = '="02/14/2024"'
$ws.Range("F1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("F1").ClearContents()

# Task 11 renamed from Login Backend to Logout Backend
$ws.Range("A11").Value = "Task 11: Logout Backend"

# Task 12 gets a real label now (was just "Task 12:")
$ws.Range("A12").Value = "Task 12 :Login Backend"

# Update column C width (target XML width 64.7109375 characters)
$ws.Columns.Item(3).ColumnWidth = 63.88

# Update selected cell
$ws.Range("C13").Select()
